# Commit message: "finished generating two results datasets"
# The v1 ("firm_properties_gemini_with_grounding_v1") row is removed from the
# VersionNotes table, leaving the v2 and v3 rows. Those two remaining rows now
# show that both results datasets finished generating (Complete? = "Y"), and
# the v3 row additionally records that it did not use Google Search ("N").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (firm_properties_gemini_with_grounding_v1).
$ws.Rows.Item(2).Delete() | Out-Null

# The former row 3 (v2) is now row 2; the former row 4 (v3) is now row 3.
# Both datasets have finished generating, so mark Complete? = "Y" for both,
# and record that the v3 row did not use Google Search grounding ("N").
$ws.Range("J2").Value = "Y"
$ws.Range("I3").Value = "N"
$ws.Range("J3").Value = "Y"

# Update the view: clear the frozen/scrolled topLeftCell and move the
# selection to A7.
$ws.Range("A7").Select() | Out-Null
